$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptos list (prices in column D, 1h volume % in column E).
# A leading apostrophe forces plain numeric-looking strings to be stored
# as text (matching the source data, which is text, not numbers) so
# Excel doesn't silently convert them to floating point numbers and lose
# formatting (trailing zeros, exact decimal value, etc.).
$ws.Range("D2").Value = '42.631.11'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '2.282.03'
$ws.Range("E3").Value = '  -3.88%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''300.58'
$ws.Range("E5").Value = '  -3.04%  '
$ws.Range("D6").Value = '''97.15'
$ws.Range("E6").Value = '  -6.69%  '
$ws.Range("D7").Value = '''0.502'
$ws.Range("E7").Value = '  -1.95%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''0.498'
$ws.Range("E9").Value = '  -4.38%  '
$ws.Range("D10").Value = '''33.60'
$ws.Range("E10").Value = '  -6.03%  '
$ws.Range("D11").Value = '''0.0788'
$ws.Range("E11").Value = '  -2.38%  '
$ws.Range("D12").Value = '''50.72'
$ws.Range("E12").Value = '  -4.78%  '
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").Value = '''6.65'
$ws.Range("E14").Value = '  -4.22%  '
$ws.Range("D15").Value = '2.636.27'
$ws.Range("E15").Value = '  -3.94%  '
$ws.Range("D16").Value = '''15.27'
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("D17").Value = '2.281.43'
$ws.Range("E17").Value = '  -3.87%  '
$ws.Range("D18").Value = '''0.787'
$ws.Range("E18").Value = '  -2.68%  '
$ws.Range("D19").Value = '42.506.25'
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").Value = '0.0₃0894'
$ws.Range("E20").Value = '  -2.02%  '
$ws.Range("D21").Value = '''11.43'
$ws.Range("E21").Value = '  -3.95%  '
$ws.Range("D22").Value = '''5.99'
$ws.Range("E22").Value = '  -4.93%  '
$ws.Range("D23").Value = '''66.65'
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").Value = '''235.56'
$ws.Range("E24").Value = '  -1.86%  '
$ws.Range("D25").Value = '''1.93'
$ws.Range("E25").Value = '  -5.26%  '
$ws.Range("E26").Value = '  -4.46%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").Value = '''24.39'
$ws.Range("E28").Value = '  -4.91%  '
$ws.Range("D29").Value = '''2.30'
$ws.Range("E29").Value = '  -1.14%  '
$ws.Range("D30").Value = '''165.03'
$ws.Range("E30").Value = '  +2.77%  '
$ws.Range("D31").Value = '''33.60'
$ws.Range("E31").Value = '  -8.03%  '
$ws.Range("E32").Value = '  -3.86%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("D34").Value = '''4.95'
$ws.Range("E34").Value = '  -5.00%  '
$ws.Range("E36").Value = '  -5.27%  '
$ws.Range("D37").Value = '''4.32'
$ws.Range("E37").Value = '  -7.27%  '
$ws.Range("E38").Value = '  -8.25%  '
$ws.Range("D39").Value = '''16.17'
$ws.Range("E39").Value = '  -11.10%  '
$ws.Range("D40").Value = '''0.0999'
$ws.Range("E40").Value = '  -5.19%  '
$ws.Range("E41").Value = '  -8.27%  '
$ws.Range("E42").Value = '  -3.28%  '
$ws.Range("E43").Value = '  -8.77%  '
$ws.Range("D44").Value = '1.958.80'
$ws.Range("E44").Value = '  -3.55%  '
$ws.Range("D45").Value = '''0.0282'
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").Value = '''17.70'
$ws.Range("E46").Value = '  -10.29%  '
$ws.Range("D47").Value = '''9.67'
$ws.Range("E47").Value = '  -8.29%  '
$ws.Range("D48").Value = '''2.83'
$ws.Range("E48").Value = '  -9.15%  '
$ws.Range("E49").Value = '  -3.36%  '

# Row 50 / 51: THORChain dropped one rank, a new row (MultiversX) was
# inserted above it, and RocketPoolETH was removed from the bottom of
# the list.
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").Value = '''53.12'
$ws.Range("E50").Value = '  -7.86%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = '''4.68'
$ws.Range("E51").Value = '  -1.28%  '
